$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the scoring columns (J:N = Lighting, Texture/Detail, Anatomy,
# Style Match, Artifacts) for the four rows that already have results.
$ws.Range("J2").Value = 4
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 3
$ws.Range("M2").Value = 3
$ws.Range("N2").Value = 3

$ws.Range("J9").Value = 2
$ws.Range("K9").Value = 3
$ws.Range("L9").Value = 3
$ws.Range("M9").Value = 3
$ws.Range("N9").Value = 2

$ws.Range("J16").Value = 3
$ws.Range("K16").Value = 4
$ws.Range("L16").Value = 3
$ws.Range("M16").Value = 3
$ws.Range("N16").Value = 3

$ws.Range("J23").Value = 1
$ws.Range("K23").Value = 1
$ws.Range("L23").Value = 3
$ws.Range("M23").Value = 3
$ws.Range("N23").Value = 2

# Widen columns H:O so the new score/notes columns are readable.
$ws.Columns.Item(8).ColumnWidth = 13
$ws.Columns.Item(9).ColumnWidth = 15.8333333333333
$ws.Columns.Item(10).ColumnWidth = 21.5
$ws.Columns.Item(11).ColumnWidth = 20.5
$ws.Columns.Item(12).ColumnWidth = 18.5
$ws.Columns.Item(13).ColumnWidth = 15.8333333333333
$ws.Columns.Item(14).ColumnWidth = 31.3333333333333
$ws.Columns.Item(15).ColumnWidth = 18.1666666666667

# Leave the cursor where the author's last edit was.
$ws.Range("N23").Select()
